$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.720.14'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '2.286.51'
$ws.Range("E3").Value = '  -0.23%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '115.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +11.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '268.46'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.84%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.632'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.27%  '
$ws.Range("E8").Value = '  +0.23%  '
$ws.Range("E9").Value = '  +2.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '49.08'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0944'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.98'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +12.71%  '
$ws.Range("E13").Value = '  +0.92%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.90'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.71%  '
$ws.Range("D15").Value = '2.632.23'
$ws.Range("E15").Value = '  -0.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.883'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.85%  '
$ws.Range("D17").Value = '2.287.16'
$ws.Range("E17").Value = '  -0.55%  '
$ws.Range("D18").Value = '43.604.51'
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("E19").Value = '  -0.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.03'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +12.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.39'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.14%  '
$ws.Range("E22").Value = '  -2.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.94'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '233.46'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.90'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.67'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.88%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '42.07'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.59%  '
$ws.Range("E29").Value = '  -1.66%  '
$ws.Range("E30").Value = '  +0.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '173.30'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.97%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.60'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.07%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0929'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.92%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.74'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.68%  '
$ws.Range("E35").Value = '  +0.65%  '
$ws.Range("E36").Value = '  -5.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0356'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.107'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.77'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '14.73'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +19.72%  '
$ws.Range("E41").Value = '  +4.66%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.86'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +14.25%  '
$ws.Range("E43").Value = '  +1.45%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.37'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +21.41%  '
$ws.Range("E45").Value = '  +0.16%  '
$ws.Range("E46").Value = '  +0.43%  '
$ws.Range("E47").Value = '  -1.21%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '103.18'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.33%  '
$ws.Range("B49").Value = 'TrustWalletToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.27'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.87%  '
$ws.Range("E50").Value = '  -1.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.457'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.72%  '
